$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows for the "2021-01-02" prediction block (rows 46-53),
# pushing the existing "2021-01-09" block down to rows 54-61.
$ws.Rows("46:53").Insert()

# Force column A to be stored as text (it looks like a date and Excel
# would otherwise silently convert "2021-01-02" into a date serial).
$ws.Range("A46:A53").NumberFormat = "@"

$ws.Range("A46").Value = "2021-01-02"
$ws.Range("B46").Value = "03 Jan -- 09 Jan 2021"
$ws.Range("C46").Value = 94.56999999999999
$ws.Range("D46").Value = 258.24
$ws.Range("E46").Value = 163.66
$ws.Range("F46").Value = "KNN"
$ws.Range("G46").Value = 1.42
$ws.Range("H46").Value = 37.08
$ws.Range("I46").Value = 45.77
$ws.Range("J46").Value = 144.14
$ws.Range("K46").Value = 143.76

$ws.Range("A47").Value = "2021-01-02"
$ws.Range("B47").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("D47").Value = 265.22
$ws.Range("F47").Value = "KNN"

$ws.Range("A48").Value = "2021-01-02"
$ws.Range("B48").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D48").Value = 267.99
$ws.Range("F48").Value = "KNN"

$ws.Range("A49").Value = "2021-01-02"
$ws.Range("B49").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D49").Value = 260.95
$ws.Range("F49").Value = "KNN"

$ws.Range("A50").Value = "2021-01-02"
$ws.Range("B50").Value = "31 Jan -- 06 Feb 2021"
$ws.Range("D50").Value = 254.08
$ws.Range("F50").Value = "KNN"

$ws.Range("A51").Value = "2021-01-02"
$ws.Range("B51").Value = "07 Feb -- 13 Feb 2021"
$ws.Range("D51").Value = 252.24
$ws.Range("F51").Value = "KNN"

$ws.Range("A52").Value = "2021-01-02"
$ws.Range("B52").Value = "14 Feb -- 20 Feb 2021"
$ws.Range("D52").Value = 250.11
$ws.Range("F52").Value = "KNN"

$ws.Range("A53").Value = "2021-01-02"
$ws.Range("B53").Value = "21 Feb -- 27 Feb 2021"
$ws.Range("D53").Value = 246.81
$ws.Range("F53").Value = "KNN"

# Restore default (General) number format/style on column A so the
# cells match the workbook default style used throughout the sheet.
$ws.Range("A46:A53").Style = "Normal"

$ws.Range("A1").Select()
